$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Withdraw"
$ws.Range("D2").Value = "1 month"

$ws.Range("D3").Select()
